$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header strings
$ws.Range("C1").Value = "Avg Haus. Dist"
$ws.Range("D1").Value = "Max"

# Update data values
$ws.Range("C3").Value = 2.8
$ws.Range("D3").Value = 20

$ws.Range("C4").Value = 2.3
$ws.Range("D4").Value = 24

$ws.Range("C5").Value = 2.4
$ws.Range("D5").Value = 17.7

$ws.Range("C6").Value = 2.9
$ws.Range("D6").Value = 18.1

# Update Row8 label from "Average" -> actually already "Average" stays same meaning but index changes; keep value
$ws.Range("A8").Value = "Average"

# Add Row9 StdDev
$ws.Range("A9").Value = "StdDev"
$ws.Range("C9").Formula = "=STDEV.S(C3:C6)"
$ws.Range("D9").Formula = "=STDEV.S(D3:D6)"
